# Fruta / hortaliza, semanal
# Insert a new weekly data point at row 215 (pushing the existing rows
# 215:286 down to 216:287) on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a new row at 215.
$ws.Rows("215:215").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A215").Value = 7
$ws.Range("B215").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C215").Value = "Ñuble"
$ws.Range("D215").Value = 44985
$ws.Range("E215").Value = 16
$ws.Range("F215").Value = 100112032
$ws.Range("G215").Value = "Zapallo italiano"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 100
$ws.Range("K215").Value = 5000
$ws.Range("L215").Value = 5500
$ws.Range("M215").Value = 5250
$ws.Range("N215").Value = '$/caja 50 unidades'
$ws.Range("O215").Value = "Región del Maule"
$ws.Range("P215").Value = 105
$ws.Range("Q215").Value = 50
$ws.Range("R215").Value = "Hortaliza"
